$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.617.43"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "1.851.68"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "264.81"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5250"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3253"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06809"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.08"
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7835"
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07800"
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").Value = "1.841.97"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.66"
$ws.Range("E14").Value = "  -0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.033"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.03"
$ws.Range("E17").Value = "  -0.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007993"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9999"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "26.647.85"
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.651"
$ws.Range("E21").Value = "  +2.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.506"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.028"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.02"
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.192"
$ws.Range("E25").Value = "  -6.33%  "
$ws.Range("E26").Value = "  +2.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.08"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "112.17"
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.204"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08739"
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04843"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7238"
$ws.Range("E33").Value = "  +4.96%  "
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.879"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.115"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.263"
$ws.Range("E37").Value = "  +2.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01798"
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.4886"
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9045"
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "111.19"
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.990"
$ws.Range("E42").Value = "  -3.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9998"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.694"
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4222"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.076"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05889"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1238"
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.13"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.8896"
$ws.Range("E50").Value = "  +3.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "60.16"
$ws.Range("E51").Value = "  +1.42%  "
